$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update incexp row (row 2/3, column B)
$ws.Range("B2").Value = "39.11***"
$ws.Range("B3").Value = "(8.47)"

# Update incvar row (row 4/5, column C)
$ws.Range("C4").Value = "1.86***"
$ws.Range("C5").Value = "(0.46)"

# Update rincvar row (row 6/7, column D)
$ws.Range("D6").Value = "2.49***"
$ws.Range("D7").Value = "(0.35)"
